$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "$ 27.489 CLP 16-11-20",
    "$ 34.434 CLP 16-11-20",
    "$ 27.489 CLP 16-11-20",
    "$ 27.495 CLP 17-11-20",
    "$ 34.442 CLP 17-11-20",
    "$ 27.495 CLP 17-11-20",
    "$ 34.458 CLP 19-11-20",
    "$ 27.508 CLP 19-11-20",
    "$ 34.458 CLP 19-11-20",
    "$ 27.508 CLP 19-11-20",
    "$ 27.508 CLP 19-11-20",
    "$ 27.515 CLP 20-11-20"
)

$startRow = 68
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
